# "Fixed bug for purchase functions"
#
# The 01-2015 issue-tracker sheet lists six "采购订单" (purchase order)
# bugs that were pending ("待解决"/blank). This commit marks them all as
# resolved ("已解决") and records the resolution date (2015-01-15).
# It also nudges the active selection on that sheet and resizes/repositions
# the first screenshot image to better match its row, plus a couple of
# row-height tweaks on the "Issue List" summary sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Issue List")
$ws2 = $wb.Worksheets.Item("01-2015")

# ---------------------------------------------------------------------
# 1. Mark the six purchase-order issues (rows 2-7) as resolved and stamp
#    the resolution date in column F (解决日期). Column E already carries
#    the record date for rows 2-3; we copy its date format into F so the
#    new value renders as m/d/yyyy instead of a raw serial number.
# ---------------------------------------------------------------------
$resolveDate = Get-Date -Year 2015 -Month 1 -Day 15 -Hour 0 -Minute 0 -Second 0

for ($r = 2; $r -le 7; $r++) {
    $ws2.Cells.Item($r, 5).Copy($ws2.Cells.Item($r, 6))
    $ws2.Cells.Item($r, 6).Value = $resolveDate
    $ws2.Cells.Item($r, 4).Value = "已解决"
}

# ---------------------------------------------------------------------
# 2. Move the active selection on the "01-2015" sheet to C7 (also drops
#    the old scrolled-down "topLeftCell" view state).
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("C7").Select()

# ---------------------------------------------------------------------
# 3. Resize/reposition the first screenshot on "01-2015" (the grid-header
#    issue picture anchored near row 2) to its new placement.
# ---------------------------------------------------------------------
$shp = $ws2.Shapes.Item(1)
$shp.Left   = 809.25
$shp.Top    = 15.75
$shp.Width  = 783.6488188976377
$shp.Height = 206.61370078740157

# ---------------------------------------------------------------------
# 4. A few row-height adjustments on the "Issue List" summary sheet.
# ---------------------------------------------------------------------
$ws1.Rows.Item(9).RowHeight  = 33
$ws1.Rows.Item(34).RowHeight = 33
$ws1.Rows.Item(36).RowHeight = 66
